$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting cand_no/last_name/first_name/phone_no right.
$ws.Columns("A:A").Insert()

# New header for the inserted "trade" column.
$ws.Range("A1").Value = "trade"

# Fill the trade column: EL for the first batch of rows (2-6), VM for the rest (7-11).
$ws.Range("A2:A6").Value = "EL"
$ws.Range("A7:A11").Value = "VM"

# The second batch's candidate numbers (now in column C) switch prefix from EL to VM.
$ws.Range("C7").Value = "VM98765"
$ws.Range("C8").Value = "VM54321"
$ws.Range("C9").Value = "VM666666"
$ws.Range("C10").Value = "VM12345"
$ws.Range("C11").Value = "VM98765"

$ws.Range("C11").Select()
